$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: age (C4) and enrollment date (D4) become real numbers instead of text ---
$ws.Range("D4").Value = 42249
$ws.Range("D4").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("C4").Value = 18
$ws.Range("C4").NumberFormat = "General"

# --- Row 5: same refactor (age / date become numbers); name/code stay text ---
$ws.Range("D5").Value = 42249
$ws.Range("D5").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("C5").Value = 18
$ws.Range("C5").NumberFormat = "General"

# --- column width for column D widened ---
$ws.Columns.Item(4).ColumnWidth = 24.4

# --- selection moves to D5 ---
$ws.Range("D5").Select() | Out-Null
